# Applies the diff to bjella_lab3.docx:
#  - 5.29 paragraph: add hanging indent, drop leading tab, replace the
#    answer text, keep a separate tab run.
#  - 5.30 paragraph: add a lastRenderedPageBreak marker before "5.30:",
#    and replace the single answer run with several split runs.
#  - 5.31 paragraph: add hanging indent, add the answer text (split from
#    the tab run), and absorb/remove the trailing empty paragraph.

$d = $word.ActiveDocument

$pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---- 5.29 ----
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*5.29:*") {
        $target = $p
        break
    }
}
$body = '<w:body><w:p><w:pPr><w:ind w:left="1440" w:hanging="720"/></w:pPr><w:r><w:t>5.29:</w:t></w:r><w:r><w:tab/></w:r><w:r><w:t>There are 600 total instruction cycles (12 per loop times 50 loops), multiplied by 0.4 microseconds per instruction cycle (4 clock cycles per instruction cycle) = 240 microseconds.</w:t></w:r></w:p></w:body>'
[void]$target.Range.InsertXML($pkgHeader + $body + $pkgFooter)

# ---- 5.30 ----
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*5.30:*") {
        $target = $p
        break
    }
}
$body = '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>5.30:</w:t></w:r><w:r><w:tab/></w:r><w:r><w:t>1.6</w:t></w:r><w:r><w:t xml:space="preserve"> microseconds, since it takes 4 fewer</w:t></w:r><w:r><w:t xml:space="preserve"> instruction</w:t></w:r><w:r><w:t xml:space="preserve"> cycles</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p></w:body>'
[void]$target.Range.InsertXML($pkgHeader + $body + $pkgFooter)

# ---- 5.31 (plus removal of the trailing empty paragraph) ----
$target = $null
$targetIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*5.31:*") {
        $target = $p
        $targetIndex = $i
        break
    }
}
$nextPara = $d.Paragraphs($targetIndex + 1)
$r = $d.Range($target.Range.Start, $nextPara.Range.End)
$body = '<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1440" w:hanging="720"/></w:pPr><w:r><w:t>5.31:</w:t></w:r><w:r><w:tab/></w:r><w:r><w:t>150 microseconds / 0.4 microseconds = 375 instruction cycles, divide by 12 and round up to 32, so set value to 0x20.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body>'
[void]$r.InsertXML($pkgHeader + $body + $pkgFooter)
